$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '86.431.44'
$ws.Range("E2").Value = '  +5.39%  '

# Row 3
$ws.Range("D3").Value = '3.275.03'
$ws.Range("E3").Value = '  +2.62%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.997'
$ws.Range("E4").Value = '  -0.39%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.83'
$ws.Range("E5").Value = '  -1.39%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '625.36'
$ws.Range("E6").Value = '  -0.15%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.377'
$ws.Range("E7").Value = '  +31.04%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.998'
$ws.Range("E8").Value = '  -0.14%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.628'
$ws.Range("E9").Value = '  +6.92%  '

# Row 10
$ws.Range("D10").Value = '3.268.23'
$ws.Range("E10").Value = '  +2.46%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.578'
$ws.Range("E11").Value = '  -2.25%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000262'
$ws.Range("E12").Value = '  +1.02%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.174'
$ws.Range("E13").Value = '  +5.13%  '

# Row 14
$ws.Range("B14").Value = 'Avalanche'
$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.35'
$ws.Range("E14").Value = '  +8.03%  '

# Row 15
$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").Value = '3.862.33'
$ws.Range("E15").Value = '  +2.10%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.21'
$ws.Range("E16").Value = '  -2.19%  '

# Row 17
$ws.Range("D17").Value = '86.216.17'
$ws.Range("E17").Value = '  +5.33%  '

# Row 18
$ws.Range("D18").Value = '3.283.40'
$ws.Range("E18").Value = '  +2.87%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.13'
$ws.Range("E19").Value = '  +0.42%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '9.14'
$ws.Range("E20").Value = '  +1.98%  '

# Row 21
$ws.Range("B21").Value = 'SuiNetwork'
$ws.Range("C21").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.91'
$ws.Range("E21").Value = '  -9.01%  '

# Row 22
$ws.Range("B22").Value = 'BitcoinCash'
$ws.Range("C22").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '433.33'
$ws.Range("E22").Value = '  -0.43%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.33'
$ws.Range("E23").Value = '  +3.92%  '

# Row 24
$ws.Range("E24").Value = '  -0.48%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.16'
$ws.Range("E25").Value = '  -2.48%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.02'
$ws.Range("E26").Value = '  +9.76%  '

# Row 27
$ws.Range("D27").Value = '3.473.31'
$ws.Range("E27").Value = '  +3.66%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '76.60'
$ws.Range("E28").Value = '  -0.15%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0000130'
$ws.Range("E29").Value = '  +5.29%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.997'
$ws.Range("E30").Value = '  -0.24%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.166'
$ws.Range("E31").Value = '  +20.42%  '

# Row 32
$ws.Range("E32").Value = '  -0.17%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '8.86'
$ws.Range("E33").Value = '  -2.05%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '544.65'
$ws.Range("E34").Value = '  -7.08%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.44'
$ws.Range("E35").Value = '  -4.32%  '

# Row 36
$ws.Range("E36").Value = '  -1.34%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.79'
$ws.Range("E37").Value = '  +10.75%  '

# Row 38
$ws.Range("E38").Value = '  -12.97%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '22.54'
$ws.Range("E39").Value = '  -1.23%  '

# Row 40
$ws.Range("B40").Value = 'WhiteBITCoin'
$ws.Range("C40").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '21.74'
$ws.Range("E40").Value = '  +4.56%  '

# Row 41
$ws.Range("B41").Value = 'FirstDigitalUSD'
$ws.Range("C41").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.991'
$ws.Range("E41").Value = '  -0.77%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.404'
$ws.Range("E42").Value = '  -1.08%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.00'
$ws.Range("E43").Value = '  -2.20%  '

# Row 44
$ws.Range("B44").Value = 'dogwifhat'
$ws.Range("C44").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.96'
$ws.Range("E44").Value = '  -3.64%  '

# Row 45
$ws.Range("B45").Value = 'USDe'
$ws.Range("C45").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.999'
$ws.Range("E45").Value = '  -0.23%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '152.53'
$ws.Range("E46").Value = '  -4.93%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '179.59'
$ws.Range("E47").Value = '  -4.53%  '

# Row 48
$ws.Range("B48").Value = 'ImmutableX'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.34'
$ws.Range("E48").Value = '  +0.39%  '

# Row 49
$ws.Range("B49").Value = 'OKB'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '44.41'
$ws.Range("E49").Value = '  -0.54%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.24'
$ws.Range("E50").Value = '  +0.52%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.628'
$ws.Range("E51").Value = '  -0.32%  '

